# pacientes.xlsx -- "tomar solo los datos del pulso" edit
#
# Drops the per-patient name column's values (keeping the column/border
# removed + General-format empty cells), flattens two previously-computed
# "running average" Pulso cells (E10/E11) back to plain numbers, and
# appends a long column of additional Pulso readings (rows 12-26) that
# the new chart will source from -- plus a couple of stray formatted-but-
# empty cells (I9, C19) that come along for the ride, matching the
# original author's manual edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Nombre column (B): blank out the patient names for the existing
#    rows and pre-format the same way for the new rows being added
#    below. Clearing the border + re-asserting "General" number format
#    in one shot on the whole range gives every cell the same new style
#    (no border, General) instead of a different style per cell.
# ---------------------------------------------------------------------
$bCol = $ws.Range("B3:B26")
$bCol.Value = $null
$bCol.Borders.LineStyle = -4142
$bCol.NumberFormat = "General"

# ---------------------------------------------------------------------
# 2. Pulso column (E): rows 10/11 were a computed running average;
#    flatten them to plain numbers and drop the "integer" number
#    format, reusing the plain bordered style already used elsewhere
#    (copy format from C10, which already carries that exact style).
# ---------------------------------------------------------------------
$ws.Range("C10").Copy()
$ws.Range("E10:E11").PasteSpecial(-4122)
$ws.Range("E10").Value = 85
$ws.Range("E11").Value = 80

# ---------------------------------------------------------------------
# 3. Append more Pulso readings in rows 12-26, taking on the same two
#    styles already present in the sheet: the "integer format + border"
#    style (copied from E5) for most rows, and the plain bordered style
#    (copied from C10) for the two rows that hold whole numbers.
# ---------------------------------------------------------------------
$ws.Range("E5").Copy()
$ws.Range("E12:E13").PasteSpecial(-4122)
$ws.Range("E16:E26").PasteSpecial(-4122)

$ws.Range("C10").Copy()
$ws.Range("E14:E15").PasteSpecial(-4122)

$ws.Range("E12").Value = 88.3333333333333
$ws.Range("E13").Value = 128.333333333333
$ws.Range("E14").Value = 85
$ws.Range("E15").Value = 80
$ws.Range("E16").Value = 88.3333333333333
$ws.Range("E17").Value = 70
$ws.Range("E18").Value = 85
$ws.Range("E19").Value = 80
$ws.Range("E20").Value = 88.3333333333333
$ws.Range("E21").Value = 93.3333333333333
$ws.Range("E22").Value = 98.3333333333333
$ws.Range("E23").Value = 103.333333333333
$ws.Range("E24").Value = 108.333333333333
$ws.Range("E25").Value = 113.333333333333
$ws.Range("E26").Value = 118.333333333333

# ---------------------------------------------------------------------
# 4. Two stray formatted-but-empty cells (underline style, no border)
#    picked up elsewhere in the sheet -- copy the format from the one
#    that already exists at L17.
# ---------------------------------------------------------------------
$ws.Range("L17").Copy()
$ws.Range("I9").PasteSpecial(-4122)
$ws.Range("C19").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 5. Column E now needs to be wide enough to show the Pulso values
#    (best-fit width, ~10.55 characters; 9.6 is the nearest input that
#    this engine's column-width quantizer resolves to that width).
# ---------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 9.6

# ---------------------------------------------------------------------
# 6. Restore the selection to where the author left it.
# ---------------------------------------------------------------------
$ws.Range("B3").Select()
